# Restore C10 value from 18 to 1 (numeric) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
